$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "...{data_inicio_pedido}. Dou o seguinte despacho:"
#      -> "...{data_inicio_pedido}, dou o seguinte despacho:"
#    i.e. the bold "}." becomes "},", and the word "Dou" is re-typed in
#    lower-case as "dou" (landing in two runs: "d" + "ou o seguinte despacho:").
# ---------------------------------------------------------------------------

$full = $d.Content.Text
$idx = $full.IndexOf("}. Dou o seguinte despacho:")
if ($idx -lt 0) { throw "anchor text '}. Dou o seguinte despacho:' not found" }

# "." -> "," (stays inside the existing bold run, merges with the trailing
# bold space run into a single "}, " run)
$rPeriod = $d.Range($idx + 1, $idx + 2)
$rPeriod.Text = ","

# Re-type "Dou" as "dou": wipe the old run's text, then insert "d" and
# "ou o seguinte despacho:" as two separate InsertAfter calls so they stay
# distinct runs (matching how Word records a retype of the capitalised word).
$full2 = $d.Content.Text
$idx2 = $full2.IndexOf("Dou o seguinte despacho")
if ($idx2 -lt 0) { throw "anchor text 'Dou o seguinte despacho' not found" }

$rWord = $d.Range($idx2, $idx2 + 24)   # "Dou o seguinte despacho:"
$rWord.Text = ""

$rInsD = $d.Range($idx2, $idx2)
$rInsD.InsertAfter("d")

$rInsRest = $d.Range($idx2 + 1, $idx2 + 1)
$rInsRest.InsertAfter("ou o seguinte despacho:")

# ---------------------------------------------------------------------------
# 2) Move the "_GoBack" bookmark from its old spot (next to "Aos
#    interessados") down to the end of the paragraph we just edited.
# ---------------------------------------------------------------------------

if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$full3 = $d.Content.Text
$endIdx = $full3.IndexOf("despacho:") + "despacho:".Length
if ($endIdx -lt 9) { throw "anchor text 'despacho:' not found" }

# A bookmark collapsed exactly on the paragraph-mark boundary can't be
# addressed directly, so park a throwaway character there first, anchor the
# bookmark just before it, then remove the throwaway character again.
$rTemp = $d.Range($endIdx, $endIdx)
$rTemp.InsertAfter("X")

$rBookmark = $d.Range($endIdx, $endIdx)
$d.Bookmarks.Add("_GoBack", $rBookmark)

$rCleanup = $d.Range($endIdx, $endIdx + 1)
$rCleanup.Text = ""
